# Scheduled market-data refresh for the Yojimbo Leve-profit tracker.
# Updates the live Universalis price columns (H:N) on each crafting-job
# sheet with freshly pulled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 36772.2
$ws.Range("I64").Value = 3194.4
$ws.Range("J64").Value = 53561.1
$ws.Range("K64").Value = 3194.4
$ws.Range("L64").Value = 53561.1
$ws.Range("M64").Value = -2946.4
$ws.Range("N64").Value = -54057.1
$ws.Range("H67").Value = 36772.2
$ws.Range("I67").Value = 3194.4
$ws.Range("J67").Value = 53561.1
$ws.Range("K67").Value = 3194.4
$ws.Range("L67").Value = 53561.1
$ws.Range("M67").Value = -2336.4
$ws.Range("N67").Value = -55277.1
$ws.Range("H94").Value = 7278.222
$ws.Range("I94").Value = 6688
$ws.Range("J94").Value = 12000
$ws.Range("K94").Value = 6688
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = -6237
$ws.Range("N94").Value = -12902
$ws.Range("H97").Value = 967.7778
$ws.Range("I97").Value = 250
$ws.Range("J97").Value = 2403.3333
$ws.Range("K97").Value = 750
$ws.Range("L97").Value = 7209.999899999999
$ws.Range("M97").Value = -254
$ws.Range("N97").Value = -8201.999899999999
$ws.Range("H115").Value = 642.5
$ws.Range("I115").Value = 642.5
$ws.Range("K115").Value = 1927.5
$ws.Range("M115").Value = -360.5
$ws.Range("H135").Value = 8249.354499999999
$ws.Range("I135").Value = 6626.4443
$ws.Range("J135").Value = 10496.462
$ws.Range("K135").Value = 59637.9987
$ws.Range("L135").Value = 94468.158
$ws.Range("M135").Value = -57102.9987
$ws.Range("N135").Value = -99538.158
$ws.Range("H137").Value = 4956.4
$ws.Range("I137").Value = 8847
$ws.Range("J137").Value = 2733.2
$ws.Range("K137").Value = 26541
$ws.Range("L137").Value = 8199.599999999999
$ws.Range("M137").Value = -23991
$ws.Range("N137").Value = -13299.6
$ws.Range("H138").Value = 9618582
$ws.Range("I138").Value = 1292.28
$ws.Range("J138").Value = 18523480
$ws.Range("K138").Value = 3876.84
$ws.Range("L138").Value = 55570440
$ws.Range("M138").Value = 1263.16
$ws.Range("N138").Value = -55580720
$ws.Range("H141").Value = 8508.549999999999
$ws.Range("I141").Value = 9269.764999999999
$ws.Range("J141").Value = 4195
$ws.Range("K141").Value = 27809.295
$ws.Range("L141").Value = 12585
$ws.Range("M141").Value = -22629.295
$ws.Range("N141").Value = -22945

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3773.88
$ws.Range("I32").Value = 3179.0674
$ws.Range("J32").Value = 8586.454
$ws.Range("K32").Value = 3179.0674
$ws.Range("L32").Value = 8586.454
$ws.Range("M32").Value = -2892.0674
$ws.Range("N32").Value = -9160.454
$ws.Range("H45").Value = 1716.3529
$ws.Range("I45").Value = 1751.9333
$ws.Range("K45").Value = 1751.9333
$ws.Range("M45").Value = -1374.9333
$ws.Range("H74").Value = 5730.7915
$ws.Range("I74").Value = 6151.3887
$ws.Range("K74").Value = 6151.3887
$ws.Range("M74").Value = -5277.3887
$ws.Range("H77").Value = 5730.7915
$ws.Range("I77").Value = 6151.3887
$ws.Range("K77").Value = 30756.9435
$ws.Range("M77").Value = -26388.9435
$ws.Range("H132").Value = 181013.81
$ws.Range("I132").Value = 4773.2563
$ws.Range("J132").Value = 562868.3
$ws.Range("K132").Value = 14319.7689
$ws.Range("L132").Value = 1688604.9
$ws.Range("M132").Value = -11789.7689
$ws.Range("N132").Value = -1693664.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 993.13513
$ws.Range("I20").Value = 957.9286
$ws.Range("J20").Value = 1102.6666
$ws.Range("K20").Value = 957.9286
$ws.Range("L20").Value = 1102.6666
$ws.Range("M20").Value = -710.9286
$ws.Range("N20").Value = -1596.6666
$ws.Range("H107").Value = 1597.3636
$ws.Range("I107").Value = 1557.1
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1557.1
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 362.9000000000001
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 50850.19
$ws.Range("I58").Value = 3890
$ws.Range("J58").Value = 201122.8
$ws.Range("K58").Value = 3890
$ws.Range("L58").Value = 201122.8
$ws.Range("M58").Value = -3687
$ws.Range("N58").Value = -201528.8
$ws.Range("H134").Value = 6481.755
$ws.Range("I134").Value = 4044.9119
$ws.Range("J134").Value = 12005.267
$ws.Range("K134").Value = 12134.7357
$ws.Range("L134").Value = 36015.801
$ws.Range("M134").Value = -9599.735700000001
$ws.Range("N134").Value = -41085.801
$ws.Range("H136").Value = 50850.19
$ws.Range("I136").Value = 3890
$ws.Range("J136").Value = 201122.8
$ws.Range("K136").Value = 11670
$ws.Range("L136").Value = 603368.3999999999
$ws.Range("M136").Value = -9120
$ws.Range("N136").Value = -608468.3999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 1490.0769
$ws.Range("I97").Value = 1530
$ws.Range("K97").Value = 1530
$ws.Range("M97").Value = -1034
$ws.Range("H122").Value = 1108.6129
$ws.Range("I122").Value = 703.2273
$ws.Range("J122").Value = 2099.5557
$ws.Range("K122").Value = 2109.6819
$ws.Range("L122").Value = 6298.6671
$ws.Range("M122").Value = 340.3181
$ws.Range("N122").Value = -11198.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 793.36365
$ws.Range("I46").Value = 679.5454999999999
$ws.Range("J46").Value = 907.1818
$ws.Range("K46").Value = 679.5454999999999
$ws.Range("L46").Value = 907.1818
$ws.Range("M46").Value = -491.5454999999999
$ws.Range("N46").Value = -1283.1818
$ws.Range("H136").Value = 2224.5908
$ws.Range("I136").Value = 1932.6316
$ws.Range("J136").Value = 2446.48
$ws.Range("K136").Value = 5797.8948
$ws.Range("L136").Value = 7339.440000000001
$ws.Range("M136").Value = -3247.8948
$ws.Range("N136").Value = -12439.44

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3442.5557
$ws.Range("I136").Value = 4178.8066
$ws.Range("J136").Value = 1812.2858
$ws.Range("K136").Value = 12536.4198
$ws.Range("L136").Value = 5436.857400000001
$ws.Range("M136").Value = -9986.4198
$ws.Range("N136").Value = -10536.8574

Write-Host "Updated market price data on ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets"
